# Remove the obsolete instrument-ID column (column A) from the
# "instruments" worksheet. The remaining columns shift left: the
# instrument-name column becomes column A, and the merged "Band ID"
# info (previously D:E) becomes C:D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete() | Out-Null

# Mirror the resulting selection left in the workbook: the whole of
# (the now former column B, now) column A is selected.
$ws.Columns.Item(1).Select() | Out-Null
